# Apply updated crypto price/volume figures per the Sun May  5 15:40:12 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.390.80'
$ws.Range('E2').Value = '  +1.08%  '

$ws.Range('D3').Value = '3.153.99'
$ws.Range('E3').Value = '  +0.75%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = "'" + '591.15'
$ws.Range('E5').Value = '  +0.29%  '

$ws.Range('D6').Value = "'" + '147.27'
$ws.Range('E6').Value = '  +0.90%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').Value = '3.151.37'
$ws.Range('E8').Value = '  +1.12%  '

$ws.Range('D9').Value = "'" + '0.531'
$ws.Range('E9').Value = '  -0.62%  '

$ws.Range('D10').Value = "'" + '0.161'
$ws.Range('E10').Value = '  -0.73%  '

$ws.Range('D11').Value = "'" + '5.96'
$ws.Range('E11').Value = '  +4.87%  '

$ws.Range('D12').Value = "'" + '0.462'
$ws.Range('E12').Value = '  -1.16%  '

$ws.Range('D13').Value = "'" + '0.0000249'
$ws.Range('E13').Value = '  -2.39%  '

$ws.Range('D14').Value = "'" + '37.21'
$ws.Range('E14').Value = '  +3.22%  '

$ws.Range('D15').Value = '3.676.45'
$ws.Range('E15').Value = '  +0.76%  '

$ws.Range('E16').Value = '  -1.09%  '

$ws.Range('E17').Value = '  +0.69%  '

$ws.Range('D18').Value = '64.149.70'
$ws.Range('E18').Value = '  +0.82%  '

$ws.Range('D19').Value = '3.154.73'
$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('D20').Value = "'" + '467.71'
$ws.Range('E20').Value = '  +0.64%  '

$ws.Range('D21').Value = "'" + '14.43'
$ws.Range('E21').Value = '  +1.55%  '

$ws.Range('D22').Value = "'" + '0.736'
$ws.Range('E22').Value = '  +0.17%  '

$ws.Range('D23').Value = "'" + '7.50'
$ws.Range('E23').Value = '  -0.41%  '

$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = "'" + '13.05'
$ws.Range('E24').Value = '  -1.74%  '

$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = "'" + '2.34'
$ws.Range('E25').Value = '  +8.44%  '

$ws.Range('D26').Value = "'" + '81.35'
$ws.Range('E26').Value = '  -1.11%  '

$ws.Range('E27').Value = '  +0.00%  '

$ws.Range('D28').Value = "'" + '9.74'
$ws.Range('E28').Value = '  +12.02%  '

$ws.Range('E29').Value = '  +0.58%  '

$ws.Range('E30').Value = '  +0.62%  '

$ws.Range('D31').Value = "'" + '7.30'
$ws.Range('E31').Value = '  +7.19%  '

$ws.Range('D32').Value = "'" + '1.00'
$ws.Range('E32').Value = '  +0.00%  '

$ws.Range('D33').Value = "'" + '27.54'
$ws.Range('E33').Value = '  +1.78%  '

$ws.Range('E34').Value = '  +0.97%  '

$ws.Range('D35').Value = '0.0₃0850'
$ws.Range('E35').Value = '  -1.18%  '

$ws.Range('D36').Value = "'" + '1.06'

$ws.Range('D37').Value = "'" + '2.34'
$ws.Range('E37').Value = '  -2.00%  '

$ws.Range('D38').Value = "'" + '6.13'
$ws.Range('E38').Value = '  +0.23%  '

$ws.Range('D39').Value = "'" + '3.28'
$ws.Range('E39').Value = '  -2.73%  '

$ws.Range('D40').Value = "'" + '51.99'
$ws.Range('E40').Value = '  +2.29%  '

$ws.Range('D41').Value = "'" + '455.43'
$ws.Range('E41').Value = '  +1.76%  '

$ws.Range('E42').Value = '  +4.57%  '

$ws.Range('E43').Value = '  +5.90%  '

$ws.Range('D44').Value = "'" + '0.0373'
$ws.Range('E44').Value = '  +0.27%  '

$ws.Range('D45').Value = '2.930.29'
$ws.Range('E45').Value = '  +0.30%  '

$ws.Range('D46').Value = "'" + '40.36'
$ws.Range('E46').Value = '  +15.41%  '

$ws.Range('E47').Value = '  -1.14%  '

$ws.Range('D48').Value = "'" + '127.83'
$ws.Range('E48').Value = '  +1.49%  '

$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('D50').Value = "'" + '2.25'
$ws.Range('E50').Value = '  +2.78%  '

$ws.Range('E51').Value = '  -0.24%  '
